$wb = $excel.ActiveWorkbook

# --- Sheet ALC (49 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 222.72728
$ws.Range("I33").Value = 207.28572
$ws.Range("J33").Value = 249.75
$ws.Range("K33").Value = 207.28572
$ws.Range("L33").Value = 249.75
$ws.Range("M33").Value = 21.71428
$ws.Range("N33").Value = -707.75
$ws.Range("H64").Value = 3769.762
$ws.Range("I64").Value = 3572.8125
$ws.Range("J64").Value = 4400
$ws.Range("K64").Value = 3572.8125
$ws.Range("L64").Value = 4400
$ws.Range("M64").Value = -3324.8125
$ws.Range("N64").Value = -4896
$ws.Range("H67").Value = 3769.762
$ws.Range("I67").Value = 3572.8125
$ws.Range("J67").Value = 4400
$ws.Range("K67").Value = 3572.8125
$ws.Range("L67").Value = 4400
$ws.Range("M67").Value = -2714.8125
$ws.Range("N67").Value = -6116
$ws.Range("H76").Value = 3120.14
$ws.Range("I76").Value = 2577.2856
$ws.Range("J76").Value = 4386.8
$ws.Range("K76").Value = 2577.2856
$ws.Range("L76").Value = 4386.8
$ws.Range("M76").Value = -2262.2856
$ws.Range("N76").Value = -5016.8
$ws.Range("H79").Value = 3120.14
$ws.Range("I79").Value = 2577.2856
$ws.Range("J79").Value = 4386.8
$ws.Range("K79").Value = 2577.2856
$ws.Range("L79").Value = 4386.8
$ws.Range("M79").Value = -1485.2856
$ws.Range("N79").Value = -6570.8
$ws.Range("H86").Value = 24720.48
$ws.Range("I86").Value = 14973.214
$ws.Range("J86").Value = 37126.09
$ws.Range("K86").Value = 14973.214
$ws.Range("L86").Value = 37126.09
$ws.Range("M86").Value = -13850.214
$ws.Range("N86").Value = -39372.09
$ws.Range("H89").Value = 24720.48
$ws.Range("I89").Value = 14973.214
$ws.Range("J89").Value = 37126.09
$ws.Range("K89").Value = 74866.07000000001
$ws.Range("L89").Value = 185630.45
$ws.Range("M89").Value = -69250.07000000001
$ws.Range("N89").Value = -196862.45

# --- Sheet ARM (26 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6345.204
$ws.Range("I32").Value = 5671.4727
$ws.Range("K32").Value = 5671.4727
$ws.Range("M32").Value = -5384.4727
$ws.Range("H63").Value = 3894.4119
$ws.Range("I63").Value = 4031.1538
$ws.Range("K63").Value = 4031.1538
$ws.Range("M63").Value = -3345.1538
$ws.Range("H66").Value = 3894.4119
$ws.Range("I66").Value = 4031.1538
$ws.Range("K66").Value = 20155.769
$ws.Range("M66").Value = -16723.769
$ws.Range("H102").Value = 2139.8647
$ws.Range("I102").Value = 2102.1714
$ws.Range("J102").Value = 2799.5
$ws.Range("K102").Value = 2102.1714
$ws.Range("L102").Value = 2799.5
$ws.Range("M102").Value = -480.1714000000002
$ws.Range("N102").Value = -6043.5
$ws.Range("H132").Value = 1451.1864
$ws.Range("I132").Value = 1099.6136
$ws.Range("J132").Value = 2482.4666
$ws.Range("K132").Value = 3298.8408
$ws.Range("L132").Value = 7447.399800000001
$ws.Range("M132").Value = -768.8407999999999
$ws.Range("N132").Value = -12507.3998

# --- Sheet BSM (15 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2679.0833
$ws.Range("I105").Value = 2687.9
$ws.Range("J105").Value = 2635
$ws.Range("K105").Value = 2687.9
$ws.Range("L105").Value = 2635
$ws.Range("M105").Value = -940.9000000000001
$ws.Range("N105").Value = -6129
$ws.Range("H113").Value = 1875
$ws.Range("I113").Value = 1875
$ws.Range("K113").Value = 1875
$ws.Range("M113").Value = 295
$ws.Range("H140").Value = 53061.668
$ws.Range("J140").Value = 53061.668
$ws.Range("L140").Value = 53061.668
$ws.Range("N140").Value = -63421.668

# --- Sheet CRP (11 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1369.2059
$ws.Range("I132").Value = 726
$ws.Range("J132").Value = 5099.8
$ws.Range("K132").Value = 2178
$ws.Range("L132").Value = 15299.4
$ws.Range("M132").Value = 352
$ws.Range("N132").Value = -20359.4
$ws.Range("H140").Value = 62738.91
$ws.Range("J140").Value = 62738.91
$ws.Range("L140").Value = 62738.91
$ws.Range("N140").Value = -73098.91

# --- Sheet GSM (33 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5112.1133
$ws.Range("I70").Value = 4400.724
$ws.Range("J70").Value = 5971.7085
$ws.Range("K70").Value = 4400.724
$ws.Range("L70").Value = 5971.7085
$ws.Range("M70").Value = -4130.724
$ws.Range("N70").Value = -6511.7085
$ws.Range("H73").Value = 5112.1133
$ws.Range("I73").Value = 4400.724
$ws.Range("J73").Value = 5971.7085
$ws.Range("K73").Value = 4400.724
$ws.Range("L73").Value = 5971.7085
$ws.Range("M73").Value = -3464.724
$ws.Range("N73").Value = -7843.7085
$ws.Range("H80").Value = 2840
$ws.Range("I80").Value = 2800
$ws.Range("K80").Value = 2800
$ws.Range("M80").Value = -1802
$ws.Range("H83").Value = 2840
$ws.Range("I83").Value = 2800
$ws.Range("K83").Value = 14000
$ws.Range("M83").Value = -9008
$ws.Range("H132").Value = 1673.2982
$ws.Range("I132").Value = 1496.6471
$ws.Range("J132").Value = 1934.4348
$ws.Range("K132").Value = 4489.9413
$ws.Range("L132").Value = 5803.3044
$ws.Range("M132").Value = -1959.9413
$ws.Range("N132").Value = -10863.3044
$ws.Range("H135").Value = 39558.57
$ws.Range("J135").Value = 39558.57
$ws.Range("L135").Value = 39558.57
$ws.Range("N135").Value = -49698.57

# --- Sheet LTW (29 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 955.8
$ws.Range("I61").Value = 980.73334
$ws.Range("J61").Value = 881
$ws.Range("K61").Value = 980.73334
$ws.Range("L61").Value = 881
$ws.Range("M61").Value = -778.73334
$ws.Range("N61").Value = -1285
$ws.Range("H100").Value = 2493.9375
$ws.Range("I100").Value = 2430.963
$ws.Range("J100").Value = 2834
$ws.Range("K100").Value = 2430.963
$ws.Range("L100").Value = 2834
$ws.Range("M100").Value = -1889.963
$ws.Range("N100").Value = -3916
$ws.Range("H113").Value = 955.8
$ws.Range("I113").Value = 980.73334
$ws.Range("J113").Value = 881
$ws.Range("K113").Value = 980.73334
$ws.Range("L113").Value = 881
$ws.Range("M113").Value = 1189.26666
$ws.Range("N113").Value = -5221
$ws.Range("H127").Value = 42651.668
$ws.Range("J127").Value = 42651.668
$ws.Range("L127").Value = 42651.668
$ws.Range("N127").Value = -52571.668
$ws.Range("H133").Value = 38684.75
$ws.Range("J133").Value = 38684.75
$ws.Range("L133").Value = 38684.75
$ws.Range("N133").Value = -43744.75

# --- Sheet WVR (8 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1542.4509
$ws.Range("I132").Value = 1054.3334
$ws.Range("K132").Value = 3163.0002
$ws.Range("M132").Value = -633.0001999999999
$ws.Range("H137").Value = 54434.125
$ws.Range("J137").Value = 54434.125
$ws.Range("L137").Value = 54434.125
$ws.Range("N137").Value = -64634.125
